# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 07:20"

# Row 25: Chequia - update Casos totales / Nuevos casos / Recuperados
$ws.Range("B25").Value = 2663
$ws.Range("C25").Value = 32
$ws.Range("E25").Value = 2641

# Rows 36/37: Tailandia overtakes Rusia in ranking, so the two rows swap
# places. Row 36 now holds Tailandia (with fresh data) and row 37 now
# holds Rusia (with the data that used to belong to row 36).
$ws.Range("A36").Value = "Tailandia"
$ws.Range("B36").Value = 1388
$ws.Range("C36").Value = 143
$ws.Range("D36").Value = 97
$ws.Range("E36").Value = 1285
$ws.Range("F36").Value = 11
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 6

$ws.Range("A37").Value = "Rusia"
$ws.Range("B37").Value = 1264
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 49
$ws.Range("E37").Value = 1211
$ws.Range("F37").Value = 8
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 4

# Row 48: Singapur - update Recuperados / Muertes hoy / Muertes
$ws.Range("E48").Value = 601
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 3
